$wb = $excel.ActiveWorkbook

# Sheet 1: ROW35-FE-LIFTER - append new log row 95
$ws1 = $wb.Worksheets.Item(1)
$ws1.Range("A95").Value = 45771.94656436342
$ws1.Range("A95").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws1.Range("B95").Value = "0x01,0x90"
$ws1.Range("C95").Value = "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x14,0x41,0x0c,"
$ws1.Range("D95").Value = "0x01,0x5a"
$ws1.Range("E95").Value = "0xd"
$ws1.Range("F95").Value = 400
$ws1.Range("G95").Value = [double]"5.68631262647114e+23"
$ws1.Range("H95").Value = 346
$ws1.Range("I95").Value = 13

# Sheet 2: ROW35-MID-LIFTER - append new log row 95
$ws2 = $wb.Worksheets.Item(2)
$ws2.Range("A95").Value = 45771.80247265047
$ws2.Range("A95").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws2.Range("B95").Value = "0x01,0x90"
$ws2.Range("C95").Value = "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x15,0x41,0x0c,"
$ws2.Range("D95").Value = "0x01,0x56"
$ws2.Range("E95").Value = "0xe"
$ws2.Range("F95").Value = 400
$ws2.Range("G95").Value = [double]"5.68631262647114e+23"
$ws2.Range("H95").Value = 342
$ws2.Range("I95").Value = 14

# Sheet 3: ROW02-FE-LIFTER - append new log row 95
$ws3 = $wb.Worksheets.Item(3)
$ws3.Range("A95").Value = 45771.94690785879
$ws3.Range("A95").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws3.Range("B95").Value = "0x01,0x90"
$ws3.Range("C95").Value = "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x0b,0x40,0x0c,"
$ws3.Range("D95").Value = "0x01,0x5a"
$ws3.Range("E95").Value = "0x3"
$ws3.Range("F95").Value = 400
$ws3.Range("G95").Value = [double]"5.68631262647114e+23"
$ws3.Range("H95").Value = 346
$ws3.Range("I95").Value = 3

# Sheet 4: ROW02-MID-LIFTER - append new log row 95
$ws4 = $wb.Worksheets.Item(4)
$ws4.Range("A95").Value = 45772.01020646991
$ws4.Range("A95").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws4.Range("B95").Value = "0x01,0x90"
$ws4.Range("C95").Value = "0xd0,0x97,0x78,0x01,0x00,0x00,0x0e,0x3f,0x0c,0x0c,"
$ws4.Range("D95").Value = "0x01,0x56"
$ws4.Range("E95").Value = "0x3"
$ws4.Range("F95").Value = 400
$ws4.Range("G95").Value = [double]"9.85046333984776e+23"
$ws4.Range("H95").Value = 342
$ws4.Range("I95").Value = 3
